# Add "Current Fiscal Year" row (row 10) to Sheet1, matching the style of
# the row above it (row 9), and select the newly-filled range -- mirroring
# what a user typing this in Excel would produce.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "Current Fiscal Year"
$ws.Range("B10").Value = 2023

# Copy row 9's formatting down into row 10 so the new label cell picks up
# the same font/alignment as the other question labels above it.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10:B10").Select() | Out-Null
